$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Fill in the previously-empty Query cells for the Login / register rows
$ws.Range("G4").Value = "User.findOne({email: email, password: hashed})"
$ws.Range("G5").Value = "user = new User(); user.save()"

# Move the usage mark for the "register" row from the "-/+" column to the "-" column
$ws.Range("J5").Value = ""
$ws.Range("I5").Value = "x"

# Column A becomes a narrow spacer column
$ws.Columns.Item(1).ColumnWidth = 4

# View tweaks
$ws.Application.ActiveWindow.Zoom = 140
$ws.Range("J5").Select()

$wb.Save()
